# "Generate Report for Handback" — refresh the localization-status report
# after a handback: update status text, handback timestamps, clear the
# stale "out of date" error notes, and resize the Status / Error Detail
# columns to fit the new text.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---- zh-cn sheet -----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Range("K2").Value = "2016-08-24 00:47:27"
$wsZh.Range("P2").Value = ""

# ---- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Range("K2").Value = "2016-08-24 00:47:34"
$wsDe.Range("P2").Value = ""

# ---- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack

# ---- Resize columns to fit the refreshed content -----------------------
# Status column widened (longer "Handed back..." text), Error Detail
# column narrowed back down now that the stale error text is gone.
$wsZh.Columns.Item(3).ColumnWidth = 29.1667
$wsZh.Columns.Item(16).ColumnWidth = 12.8333

$wsDe.Columns.Item(3).ColumnWidth = 29.1667
$wsDe.Columns.Item(16).ColumnWidth = 12.8333

$wsOverview.Columns.Item(5).ColumnWidth = 29.1667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1667
